$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new meeting-diary entry as row 9, matching the formatting of the
# previous entry (row 7): date format on A, time format on B/C, plain text
# on D, wrapped text on E.
$ws.Range("A7:E7").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A9").Value = "9/29/2023"
$ws.Range("B9").Value = 0.45833333333333331   # 11:00 AM
$ws.Range("C9").Value = 0.5                   # 12:00 PM
$ws.Range("D9").Value = "All"
$ws.Range("E9").Value = "Perform merging `ndivide workload for IDA part and set deadline"

$ws.Rows(9).RowHeight = 51

# Move the active selection the way the author left it after typing the row.
$ws.Range("E10").Select()
